# Update the value in B3 (NIM) and C3 (Nama Lengkap) for the second data row,
# and move the active selection from D2 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 data: change NIM and Name
$ws.Range("B3").Value = 2341760196
$ws.Range("C3").Value = "Kemal S"

# Update the selected/active cell shown in the saved view
$ws.Range("B2").Select()
